$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three coupon-validity rows (A2:A4) that are no longer used
$ws.Range("A2:A4").ClearContents()

# Move the active selection to A5, matching the saved view state
$ws.Range("A5").Select()
